$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 10.358308170675899
$ws.Range("B2").Value = 10.817693606886399
$ws.Range("C2").Value = 10.805419237911201
$ws.Range("A3").Value = 10.636972111434799
$ws.Range("B3").Value = 10.3703757655858
$ws.Range("C3").Value = 10.1987171422667
$ws.Range("A4").Value = 10.5438610050881
$ws.Range("B4").Value = 10.8220687868012
$ws.Range("C4").Value = 10.421240876018899
$ws.Range("A5").Value = 10.7710545810556
$ws.Range("B5").Value = 10.2615751879238
$ws.Range("C5").Value = 10.0743237618112
$ws.Range("A6").Value = 10.474565413259301
$ws.Range("B6").Value = 10.9283844284911
$ws.Range("C6").Value = 10.397241305355299
$ws.Range("A7").Value = 10.1780799205546
$ws.Range("B7").Value = 10.079247163372401
$ws.Range("C7").Value = 10.301505593533101
$ws.Range("A8").Value = 10.789261506812
$ws.Range("B8").Value = 10.454925880307
$ws.Range("C8").Value = 10.897990408932101
$ws.Range("A9").Value = 10.021917581472801
$ws.Range("B9").Value = 10.553036799540999
$ws.Range("C9").Value = 10.734879238782799
$ws.Range("A10").Value = 10.6633917093934
$ws.Range("B10").Value = 10.3217857608457
$ws.Range("C10").Value = 10.8173631726515
$ws.Range("A11").Value = 10.5677623480134
$ws.Range("B11").Value = 10.759644968648299
$ws.Range("C11").Value = 10.799596727032901
$ws.Range("A12").Value = 10.165641863347
$ws.Range("B12").Value = 10.271288389806299
$ws.Range("C12").Value = 10.9285038520688
$ws.Range("A13").Value = 10.666901345851199
$ws.Range("B13").Value = 10.745501621902999
$ws.Range("C13").Value = 10.8587942992372
$ws.Range("A14").Value = 10.969768694676899
$ws.Range("B14").Value = 10.0370370585053
$ws.Range("C14").Value = 10.379934908432
$ws.Range("A15").Value = 10.3606630952758
$ws.Range("B15").Value = 10.8586974590709
$ws.Range("C15").Value = 10.7510616211956
$ws.Range("A16").Value = 10.9493368693512
$ws.Range("B16").Value = 10.3081493099084
$ws.Range("C16").Value = 10.3142801955322
$ws.Range("A17").Value = 10.5972593135164
$ws.Range("B17").Value = 10.507003369515299
$ws.Range("C17").Value = 10.5965413126671
$ws.Range("A18").Value = 10.3712652505412
$ws.Range("B18").Value = 10.171090887223301
$ws.Range("C18").Value = 10.2065031298085
$ws.Range("A19").Value = 10.3387879713597
$ws.Range("B19").Value = 10.0205482535294
$ws.Range("C19").Value = 10.6919199448041
$ws.Range("A20").Value = 10.9675532082861
$ws.Range("B20").Value = 10.2070674347018
$ws.Range("C20").Value = 10.299232624862899
$ws.Range("A21").Value = 10.933220614607899
$ws.Range("B21").Value = 10.0226814691469
$ws.Range("C21").Value = 10.8332020729219
$ws.Range("A22").Value = 10.048647202669001
$ws.Range("B22").Value = 10.416435849485399
$ws.Range("C22").Value = 10.7636094852928
$ws.Range("A23").Value = 10.1083348543719
$ws.Range("B23").Value = 10.3251737638432
$ws.Range("C23").Value = 10.9303447350204
$ws.Range("A24").Value = 10.718286324818701
$ws.Range("B24").Value = 10.235778875870899
$ws.Range("C24").Value = 10.8518183453593
$ws.Range("A25").Value = 10.8599877081776
$ws.Range("B25").Value = 10.8910371762748
$ws.Range("C25").Value = 10.049231060279499
$ws.Range("A26").Value = 10.333181174099201
$ws.Range("B26").Value = 10.8069498973392
$ws.Range("C26").Value = 10.555089349272199
$ws.Range("A27").Value = 10.158007673095801
$ws.Range("B27").Value = 10.145585637888701
$ws.Range("C27").Value = 10.102995388598201
$ws.Range("A28").Value = 10.191995844269
$ws.Range("B28").Value = 10.6640868295374
$ws.Range("C28").Value = 10.1270463638145
$ws.Range("A29").Value = 10.2621282518169
$ws.Range("B29").Value = 10.9515624029387
$ws.Range("C29").Value = 10.994076174852401
$ws.Range("A30").Value = 10.617467737147001
$ws.Range("B30").Value = 10.630476657684699
$ws.Range("C30").Value = 10.109181969410299
$ws.Range("A31").Value = 10.582719788178601
$ws.Range("B31").Value = 10.2135319014301
$ws.Range("C31").Value = 10.7650016843701

$ws.Range("A2:C31").Select() | Out-Null
